$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the entries that are no longer needed from the planet list,
# leaving the rows blank (matching the rows that used to hold
# "Antartic Planet", "Industrial - Winter Planet", "Earth-Gold Planet",
# "Gold like Industrial Planet" and "Metalic Gold Planet").
$ws.Range("A3:B3").ClearContents()
$ws.Range("A9:B9").ClearContents()
$ws.Range("A13:B13").ClearContents()
$ws.Range("A17:B17").ClearContents()
$ws.Range("A18:B18").ClearContents()

# Update the active cell selection to match the saved state.
$ws.Range("E15").Select()
